$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the value in E6 (clear its content)
$ws.Range("E6").ClearContents()

# Update the selection to D17 as shown in the diff
$ws.Range("D17").Select()
